$d = $word.ActiveDocument

# --- 1) Relocate the "_GoBack" bookmark ---------------------------------
# In the edited document, "_GoBack" (Word's "last edit position" bookmark)
# moved from the end of the paragraph to wrap the stretch of text that was
# last touched: right before the "Y" of "Yellow" through right before the
# "B" of the second "Honey Bees".  Re-adding a bookmark named "_GoBack"
# removes the old one and plants the new one, exactly like Word itself does.

$bmStartRange = $d.Content
$bmStartRange.Find.Execute("Red, Y", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmStart = $bmStartRange.End - 1

$bmEndRange = $d.Content
$bmEndRange.Find.Execute("Honey B", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmEnd = $bmEndRange.End - 1

$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmEnd))

# --- 2) Append the new sentences about weed control at the paragraph end -
$tail = $d.Content
$tail.Find.Execute("of a teaspoon of honey. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Collapse(0)

$newRuns = @(
  "Boiling water will kill weeds just effectively as chemicals. There are also lots of other sprays that are non-chemicals like",
  " garlic spray, Oil",
  " spray,",
  " Milk spray, natural ant deterrent, and herbic",
  "i",
  "de ",
  "alternatives for weeds."
)

foreach ($run in $newRuns) {
  $tail.InsertAfter($run)
  $tail.Collapse(0)
}
